$d = $word.ActiveDocument

# 1. Header contact line: "San Jose, California" -> "San Jose, CA"
$d.Content.Find.Execute(
    "San Jose, California",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "San Jose, CA",
    2)

# 2a. Kantar bullet 1: insert line break before "Apache Kafka"
$d.Content.Find.Execute(
    "integrating real-time data streaming with Apache Kafka",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "integrating real-time data streaming with^lApache Kafka",
    2)

# 2b. Kantar bullet 2: insert line break before "integrated data quality checks"
$d.Content.Find.Execute(
    "conversion rates, and integrated data quality checks",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "conversion rates, and^lintegrated data quality checks",
    2)

# 2c. Kantar bullet 3: insert line break before "testing of individual"
$d.Content.Find.Execute(
    "error handling, and unit testing of individual pipeline components",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "error handling, and unit^ltesting of individual pipeline components",
    2)

# 2d. Kantar bullet 4: insert line break before "Bank Brazil"
$d.Content.Find.Execute(
    "ING Bank Australia, Itau Bank Brazil, Pandora UK",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "ING Bank Australia, Itau^lBank Brazil, Pandora UK",
    2)

# 3. "Data Engineer | The Sparks Foundation" -> "Data Engineer | Sparks Foundation"
$d.Content.Find.Execute(
    "Data Engineer | The Sparks Foundation",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Data Engineer | Sparks Foundation",
    2)

# 4a. Sparks Foundation bullet 1: insert line break before "improving efficiency"
$d.Content.Find.Execute(
    "NOAA and NASA datasets, improving efficiency by 25%",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "NOAA and NASA datasets,^limproving efficiency by 25%",
    2)

# 4b. Sparks Foundation bullet 2: insert line break before "reduction" AND drop "30% "
$d.Content.Find.Execute(
    "resulting in a 30% reduction in data retrieval times",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "resulting in a^lreduction in data retrieval times",
    2)

# 4c. Sparks Foundation bullet 3: insert line break before "a Star Schema"
$d.Content.Find.Execute(
    "stakeholders, including the implementation of a Star Schema",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "stakeholders, including the implementation of^la Star Schema",
    2)

# 4d. Sparks Foundation bullet 4: insert line break before "system for NASA"
$d.Content.Find.Execute(
    "developed a recommendation system for NASA Earth Observation",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "developed a recommendation^lsystem for NASA Earth Observation",
    2)

# 4e. Sparks Foundation bullet 5: insert line break before "Dockerized deployment"
$d.Content.Find.Execute(
    "including automated unit tests, Dockerized deployment",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "including automated unit tests,^lDockerized deployment",
    2)

# 5. Skills line: full rewrite
$d.Content.Find.Execute(
    "Python | SQL/NoSQL | Java | R | C | C++ | Bash | JavaScript | Scikit-learn | Keras | PyTorch | Delta Lake | Apache Iceberg | MLOps | AWS | S3 | Redshift | RDS | GCP | Azure | Snowflake | BigQuery | Databricks | dbt | Apache Hadoop | HDFS | Hive | Sqoop | HBase | MySQL | PostgreSQL | MongoDB | Elastisearch | Apache Airflow | Kafka | Spark | PySpark | Terraform | Docker | Kubernetes | Jenkins | Prometheus | Grafana | Presto | Flume | Dask | CI/CD Pipelines | Azure Blob Storage | GraphQL | Power BI",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Python | PySpark | Airflow | Kafka | C++ | Java | JavaScript | API Design and Integration | Data Warehousing | AWS | Azure | GCP | SQL/NoSQL | PostgreSQL | MongoDB | Elastisearch | Spark | Hadoop | Bash | PyTorch | Scikit-learn | Keras | Delta Lake | MLOps | Big Data Tools | BigQuery | Snowflake | dbt | HDFS | Hive | Sqoop | Hbase | Terraform | Docker | Kubernetes | Jenkins | Prometheus | Grafana | Presto | Flume | Dask | CI/CD Pipelines | Power BI | GraphQL",
    2)

# 6. Cab service bullet: insert line break before "HDFS and AWS RDS"
$d.Content.Find.Execute(
    "ingesting data into Hadoop HDFS and AWS RDS via Sqoop",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "ingesting data into Hadoop^lHDFS and AWS RDS via Sqoop",
    2)

# 7. Spar Nord Bank bullet: insert line break before "analysis to optimize"
$d.Content.Find.Execute(
    "MySQL RDS to Amazon Redshift and performed data analysis to optimize ATM refill processes",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "MySQL RDS to Amazon Redshift and performed data^lanalysis to optimize ATM refill processes",
    2)

# 8. PG Diploma date: "Aug 2020 - Sep 2021" -> "Oct 2020 - Sep 2021"
$d.Content.Find.Execute(
    "Aug 2020 - Sep 2021",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Oct 2020 - Sep 2021",
    2)

# 9. Bachelors date: "Aug 2016 - Aug 2020" -> "Aug 2016 - Jul 2020"
$d.Content.Find.Execute(
    "Aug 2016 - Aug 2020",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Aug 2016 - Jul 2020",
    2)
